# Auto-generated edit script: apply Sheets data refresh per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1729.6
$ws.Range("I70").Value = 1766.3334
$ws.Range("J70").Value = 1674.5
$ws.Range("K70").Value = 5299.0002
$ws.Range("L70").Value = 5023.5
$ws.Range("M70").Value = -5029.0002
$ws.Range("N70").Value = -5563.5

$ws.Range("H73").Value = 1729.6
$ws.Range("I73").Value = 1766.3334
$ws.Range("J73").Value = 1674.5
$ws.Range("K73").Value = 5299.0002
$ws.Range("L73").Value = 5023.5
$ws.Range("M73").Value = -4363.0002
$ws.Range("N73").Value = -6895.5

$ws.Range("H100").Value = 2916.4443
$ws.Range("I100").Value = 2853.889
$ws.Range("J100").Value = 2979
$ws.Range("K100").Value = 2853.889
$ws.Range("L100").Value = 2979
$ws.Range("M100").Value = -2312.889
$ws.Range("N100").Value = -4061

$ws.Range("H113").Value = 2166.6667
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -9508

$ws.Range("H137").Value = 3104.4146
$ws.Range("I137").Value = 2140.9412
$ws.Range("J137").Value = 7784.143
$ws.Range("K137").Value = 6422.823600000001
$ws.Range("L137").Value = 23352.429
$ws.Range("M137").Value = -3872.823600000001
$ws.Range("N137").Value = -28452.429

$ws.Range("H138").Value = 2378.4243
$ws.Range("I138").Value = 2928.2856
$ws.Range("J138").Value = 2287.859
$ws.Range("K138").Value = 8784.856800000001
$ws.Range("L138").Value = 6863.576999999999
$ws.Range("M138").Value = -3644.856800000001
$ws.Range("N138").Value = -17143.577

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3528.1738
$ws.Range("I61").Value = 3499.625
$ws.Range("J61").Value = 3543.4
$ws.Range("K61").Value = 3499.625
$ws.Range("L61").Value = 3543.4
$ws.Range("M61").Value = -3287.625
$ws.Range("N61").Value = -3967.4

$ws.Range("H122").Value = 73379.36
$ws.Range("I122").Value = 92510.09
$ws.Range("J122").Value = 3233.3333
$ws.Range("K122").Value = 277530.27
$ws.Range("L122").Value = 9699.999899999999
$ws.Range("M122").Value = -275080.27
$ws.Range("N122").Value = -14599.9999

$ws.Range("H132").Value = 5508.8276
$ws.Range("I132").Value = 6156.222
$ws.Range("J132").Value = 5217.5
$ws.Range("K132").Value = 18468.666
$ws.Range("L132").Value = 15652.5
$ws.Range("M132").Value = -15938.666
$ws.Range("N132").Value = -20712.5

$ws.Range("H136").Value = 3528.1738
$ws.Range("I136").Value = 3499.625
$ws.Range("J136").Value = 3543.4
$ws.Range("K136").Value = 10498.875
$ws.Range("L136").Value = 10630.2
$ws.Range("M136").Value = -7948.875
$ws.Range("N136").Value = -15730.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 12503258
$ws.Range("I105").Value = 15628185
$ws.Range("K105").Value = 15628185
$ws.Range("M105").Value = -15626438

$ws.Range("H134").Value = 2912.6785
$ws.Range("I134").Value = 3195.6155
$ws.Range("K134").Value = 9586.8465
$ws.Range("M134").Value = -7051.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 742.7143
$ws.Range("I16").Value = 779.8
$ws.Range("J16").Value = 650
$ws.Range("K16").Value = 779.8
$ws.Range("L16").Value = 650
$ws.Range("M16").Value = -492.8
$ws.Range("N16").Value = -1224

$ws.Range("H99").Value = 1875.5
$ws.Range("I99").Value = 1470.6666
$ws.Range("K99").Value = 1470.6666
$ws.Range("M99").Value = 27.33339999999998

$ws.Range("H113").Value = 742.7143
$ws.Range("I113").Value = 779.8
$ws.Range("J113").Value = 650
$ws.Range("K113").Value = 779.8
$ws.Range("L113").Value = 650
$ws.Range("M113").Value = 1390.2
$ws.Range("N113").Value = -4990

$ws.Range("H126").Value = 1875.5
$ws.Range("I126").Value = 1470.6666
$ws.Range("K126").Value = 4411.9998
$ws.Range("M126").Value = -1941.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 220.875
$ws.Range("I40").Value = 220.875
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 883.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -814.5
$ws.Range("N40").ClearContents()

$ws.Range("H41").Value = 2608.4285
$ws.Range("J41").Value = 2993.1667
$ws.Range("L41").Value = 8979.500100000001
$ws.Range("N41").Value = -9655.500100000001

$ws.Range("H132").Value = 2200.8518
$ws.Range("J132").Value = 2254.111
$ws.Range("L132").Value = 20286.999
$ws.Range("N132").Value = -25346.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2022.1333

$ws.Range("H132").Value = 2642.75
$ws.Range("I132").Value = 2455.0908
$ws.Range("J132").Value = 2801.5386
$ws.Range("K132").Value = 7365.2724
$ws.Range("L132").Value = 8404.6158
$ws.Range("M132").Value = -4835.2724
$ws.Range("N132").Value = -13464.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 569.0741
$ws.Range("I55").Value = 340.7
$ws.Range("J55").Value = 703.41174
$ws.Range("K55").Value = 340.7
$ws.Range("L55").Value = 703.41174
$ws.Range("M55").Value = -167.7
$ws.Range("N55").Value = -1049.41174

$ws.Range("H93").Value = 21397.2
$ws.Range("I93").Value = 21397.2
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 21397.2
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -20149.2
$ws.Range("N93").ClearContents()

$ws.Range("H100").Value = 2388.4
$ws.Range("I100").Value = 2147.5
$ws.Range("J100").Value = 2749.75
$ws.Range("K100").Value = 2147.5
$ws.Range("L100").Value = 2749.75
$ws.Range("M100").Value = -1606.5
$ws.Range("N100").Value = -3831.75

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 17148
$ws.Range("J101").Value = 17148
$ws.Range("L101").Value = 17148
$ws.Range("N101").Value = -23638

$ws.Range("H122").Value = 4336
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4670
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 14010
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -18910
